$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Update the raw counts in row 2 (VP, VN, FP, FN) with the final values.
# B2 loses its centered-cell style (reset to the default "Normal" style)
# while the value is rewritten, matching the authored edit.
$ws.Range("B2").Style = "Normal"
$ws.Range("B2").Value = 1180
$ws.Range("C2").Value = 1176
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = 2

# Re-enter the derived-metric formulas across the full column ranges so the
# unchanged rows (3:12) pick up the same shared-formula group as row 2.
$ws.Range("F2:F12").Formula = "=B2/(B2+E2)"
$ws.Range("G2:G12").Formula = "=C2/(C2+D2)"
$ws.Range("H2:H12").Formula = "=(B2+C2)/(B2+C2+D2+E2)"

# Move the sheet's active selection to match the author's final cursor spot.
$ws.Range("J18").Select()
